$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column Q (raname): rewrite English display names in the same order the
# replacement strings are first introduced (dependency order)
$ws.Range("Q2").Value = "Michael Walton"
$ws.Range("Q3:Q18").Value = "Samuel Lester"
$ws.Range("Q19").Value = "Jeffrey Campos"
$ws.Range("Q20:Q22").Value = "Christopher Garrett"
$ws.Range("Q23:Q24").Value = "Melissa Welch"
$ws.Range("Q25:Q29").Value = "Jacob Hutchinson"
$ws.Range("Q30:Q36").Value = "Alicia Graham"
$ws.Range("Q37:Q39").Value = "Kaitlin Williamson"
$ws.Range("Q40:Q41").Value = "Andrew Roman"
$ws.Range("Q42").Value = "Scott Orozco"
$ws.Range("Q43:Q46").Value = "Joseph Lee"
$ws.Range("Q47").Value = "Eileen Hill"
$ws.Range("Q48:Q55").Value = "Victor Fitzgerald"
$ws.Range("Q56:Q58").Value = "Raymond Fuller"
$ws.Range("Q59:Q68").Value = "April Flores"
$ws.Range("Q69:Q73").Value = "Robert Wright"
$ws.Range("Q74:Q87").Value = "James Becker"
$ws.Range("Q88:Q98").Value = "Travis Villarreal"
$ws.Range("Q99:Q105").Value = "Dylan Schultz"
$ws.Range("Q106:Q108").Value = "Adam Bennett"
$ws.Range("Q109:Q124").Value = "Erica Thomas"
$ws.Range("Q125:Q129").Value = "Eric Cobb"
$ws.Range("Q130:Q141").Value = "Michael Briggs"
$ws.Range("Q142:Q143").Value = "Richard Clark"
$ws.Range("Q144").Value = "Michael Taylor"
$ws.Range("Q145:Q147").Value = "Jill Travis"
$ws.Range("Q148:Q150").Value = "Breanna Jensen"
$ws.Range("Q151:Q161").Value = "Karen Bowman"
$ws.Range("Q162:Q166").Value = "Karen Webb"
$ws.Range("Q167:Q171").Value = "Joshua Cox"
$ws.Range("Q172:Q174").Value = "Sandra Hobbs"
$ws.Range("Q175:Q187").Value = "Kevin Lee"
$ws.Range("Q188:Q194").Value = "Deanna Blair"
$ws.Range("Q195").Value = "Michelle Williams"
$ws.Range("Q196:Q201").Value = "Anthony Jones"
$ws.Range("Q202:Q203").Value = "David Jensen"
$ws.Range("Q204:Q209").Value = "Kenneth Mullen"
$ws.Range("Q210").Value = "Matthew Robinson"
$ws.Range("Q211:Q217").Value = "Ashley Hopkins"
$ws.Range("Q218:Q222").Value = "Linda Smith"
$ws.Range("Q223:Q227").Value = "Aaron Nichols"
$ws.Range("Q228:Q245").Value = "Lauren Tyler"
$ws.Range("Q246").Value = "Jeffrey Carlson"
$ws.Range("Q247:Q258").Value = "Andre Howard"
$ws.Range("Q259:Q267").Value = "Stephen Suarez"
$ws.Range("Q268:Q275").Value = "Andrew Sullivan"
$ws.Range("Q276:Q282").Value = "Maria Meza"
$ws.Range("Q283:Q290").Value = "Roberta Jenkins"
$ws.Range("Q291:Q292").Value = "Meghan Dunn"
$ws.Range("Q293:Q299").Value = "Desiree Brock"
$ws.Range("Q300:Q307").Value = "Laura Watson"
$ws.Range("Q308:Q319").Value = "Valerie Cohen"
$ws.Range("Q320:Q326").Value = "Kimberly Christensen"
$ws.Range("Q327:Q334").Value = "Brian Jones"
$ws.Range("Q335:Q339").Value = "Monica Olsen"
$ws.Range("Q340:Q351").Value = "Caitlin Flores"
$ws.Range("Q352:Q358").Value = "Wendy Waters"
$ws.Range("Q359:Q362").Value = "Kristina Torres"
$ws.Range("Q363:Q367").Value = "Crystal Nielsen"
$ws.Range("Q368:Q382").Value = "Sarah Jordan"
$ws.Range("Q383:Q384").Value = "Austin Solomon"
$ws.Range("Q385:Q389").Value = "Juan Rodriguez"
$ws.Range("Q390:Q391").Value = "Jennifer Gibson"
$ws.Range("Q392:Q396").Value = "Linda Miller"
$ws.Range("Q397:Q398").Value = "Victor Martinez"
$ws.Range("Q399:Q400").Value = "Andrea Bryan"
$ws.Range("Q401:Q410").Value = "Kerry Day"
$ws.Range("Q411").Value = "Jesse Flores"
$ws.Range("Q412:Q418").Value = "Michael Grant"
$ws.Range("Q419").Value = "Michael Simmons"
$ws.Range("Q420:Q424").Value = "Scott Morales"
$ws.Range("Q425").Value = "Jose Medina"
$ws.Range("Q426:Q436").Value = "Nina Murphy"
$ws.Range("Q437:Q438").Value = "Jacob Hutchinson"
$ws.Range("Q439:Q442").Value = "Brittany Nelson"
$ws.Range("Q443:Q452").Value = "Dr. Patricia Gill"
$ws.Range("Q453").Value = "Sabrina Patterson"
$ws.Range("Q454:Q455").Value = "Morgan Hernandez"
$ws.Range("Q456:Q458").Value = "Whitney Powell"
$ws.Range("Q459:Q463").Value = "Thomas Spencer"
$ws.Range("Q464:Q473").Value = "Garrett Williams"
$ws.Range("Q474:Q477").Value = "Jill Jackson"
$ws.Range("Q478:Q484").Value = "Bradley Sullivan"
$ws.Range("Q485:Q493").Value = "Jessica Moore"
$ws.Range("Q494:Q499").Value = "Mary Hernandez"
$ws.Range("Q500:Q509").Value = "Christopher Castaneda"
$ws.Range("Q510").Value = "Danielle Wilkins"
$ws.Range("Q511:Q516").Value = "Cindy Pierce"
$ws.Range("Q517").Value = "Steve Cooper"
$ws.Range("Q518:Q521").Value = "Jessica Bradshaw"
$ws.Range("Q522:Q528").Value = "Christian Richardson"
$ws.Range("Q529:Q532").Value = "Dr. Norma Ramirez MD"
$ws.Range("Q533:Q536").Value = "Dominique Valdez"
$ws.Range("Q537:Q539").Value = "Mallory Logan"
$ws.Range("Q540").Value = "Robert Mills"
$ws.Range("Q541:Q551").Value = "Gina Rios"
$ws.Range("Q552:Q558").Value = "Linda Smith"
$ws.Range("Q559").Value = "Crystal Lopez"
$ws.Range("Q560:Q564").Value = "Deanna Cabrera"
$ws.Range("Q565:Q568").Value = "Amanda Campbell"
$ws.Range("Q569:Q570").Value = "Mary Reyes"
$ws.Range("Q571:Q583").Value = "Matthew Doyle"
$ws.Range("Q584:Q593").Value = "Julie Ewing"
$ws.Range("Q594:Q598").Value = "Deborah Lang"
$ws.Range("Q599:Q603").Value = "Brian Spears"
$ws.Range("Q604:Q605").Value = "Ryan Bennett"
$ws.Range("Q606:Q613").Value = "Valerie Sanders"
$ws.Range("Q614:Q620").Value = "Justin Spencer"
$ws.Range("Q621:Q628").Value = "Lisa Mcbride"
$ws.Range("Q629:Q632").Value = "Kenneth Owens"
$ws.Range("Q633:Q636").Value = "Victoria Malone"
$ws.Range("Q637:Q639").Value = "Victor Martinez"
$ws.Range("Q640").Value = "Aaron Hawkins"
$ws.Range("Q641:Q645").Value = "Charles Harris"
$ws.Range("Q646:Q649").Value = "Kyle Conway"
$ws.Range("Q650:Q652").Value = "Michelle Davis"
$ws.Range("Q653").Value = "Sean Russell"
$ws.Range("Q654:Q656").Value = "Leslie Callahan"
$ws.Range("Q657:Q664").Value = "Monica Wallace"
$ws.Range("Q665:Q669").Value = "Bruce English"
$ws.Range("Q670:Q671").Value = "April Dawson"
$ws.Range("Q672:Q676").Value = "Kimberly Berger"
$ws.Range("Q677:Q687").Value = "George Banks"
$ws.Range("Q688:Q692").Value = "Sheila Mendoza"
$ws.Range("Q693:Q696").Value = "Isabella Johnson"
$ws.Range("Q697:Q701").Value = "Kathleen Gonzales"
$ws.Range("Q702:Q707").Value = "Gabriela Jackson"
$ws.Range("Q708:Q712").Value = "Victoria Frederick"
$ws.Range("Q713:Q719").Value = "Krystal Kerr"
$ws.Range("Q720:Q735").Value = "Angela Velez"
$ws.Range("Q736:Q745").Value = "Edward Conway"
$ws.Range("Q746:Q749").Value = "Kristine Smith"
$ws.Range("Q750:Q758").Value = "Katherine Cole"
$ws.Range("Q759:Q767").Value = "David Wang"
$ws.Range("Q768").Value = "Paul Walter"
$ws.Range("Q769:Q770").Value = "Brenda Thompson"
$ws.Range("Q771:Q773").Value = "Ethan Tucker"
$ws.Range("Q774:Q783").Value = "David Jackson"
$ws.Range("Q784:Q787").Value = "Denise Carlson"
$ws.Range("Q788:Q792").Value = "Christopher Flores Jr."
$ws.Range("Q793:Q795").Value = "Paula Hanson"

# Column S (teamid): rewrite team UUIDs in the same order the replacement
# strings are first introduced (dependency order)
$ws.Range("S2").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S3:S4").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S5").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S6:S12").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S13").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S14").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S15").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S16:S18").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S19").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S20:S24").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S25:S27").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S28:S29").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S30").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S31:S33").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S34:S36").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S37").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S38").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S39").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S40:S41").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S42").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S43").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S44").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S45:S46").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S47").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S48:S55").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S56:S58").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S59:S68").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S69:S71").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S72").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S73").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S74:S76").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S77:S87").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S88:S90").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S91:S93").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S94:S98").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S99:S100").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S101:S102").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S103").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S104:S105").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S106:S108").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S109:S124").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S125").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S126:S127").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S128:S129").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S130").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S131").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S132:S133").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S134:S138").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S139:S141").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S142:S143").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S144").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S145").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S146").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S147").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S148:S153").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S154:S161").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S162:S166").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S167:S168").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S169").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S170").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S171").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S172:S174").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S175:S177").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S178:S184").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S185:S187").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S188").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S189:S191").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S192:S193").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S194").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S195").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S196").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S197:S201").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S202:S203").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S204:S209").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S210").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S211:S213").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S214").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S215:S217").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S218:S219").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S220").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S221:S222").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S223:S225").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S226:S227").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S228:S229").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S230:S245").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S246:S247").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S248").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S249").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S250:S252").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S253:S258").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S259:S267").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S268:S270").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S271:S273").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S274:S275").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S276").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S277").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S278:S280").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S281:S282").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S283:S288").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S289:S290").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S291:S294").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S295:S299").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S300:S306").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S307").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S308:S319").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S320:S321").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S322:S326").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S327:S339").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S340:S341").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S342").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S343").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S344").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S345").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S346:S351").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S352").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S353:S354").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S355:S358").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S359:S362").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S363:S364").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S365:S366").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S367").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S368:S370").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S371:S375").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S376:S378").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S379:S382").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S383").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S384").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S385").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S386:S389").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S390:S391").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S392").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S393:S396").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S397:S398").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S399:S400").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S401:S403").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S404:S405").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S406").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S407").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S408:S410").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S411").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S412:S414").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S415:S416").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S417:S418").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S419").Value = "de503c24-f17d-47a9-9a47-6f0a194f8c9c"
$ws.Range("S420:S421").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S422:S424").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S425").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S426:S427").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S428").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S429").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S430:S432").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S433:S436").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S437").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S438").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S439:S442").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S443:S446").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S447:S449").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S450").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S451").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S452").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S453").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S454:S455").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S456").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S457:S458").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S459:S463").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S464:S465").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S466").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S467").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S468:S469").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S470").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S471:S473").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S474:S476").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S477").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S478").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S479:S481").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S482:S484").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S485").Value = "718c6b8f-7c00-4bcb-b53c-8f3f42154362"
$ws.Range("S486:S487").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S488").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S489:S490").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S491:S493").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S494:S509").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S510").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S511:S514").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S515").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S516").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S517").Value = "718c6b8f-7c00-4bcb-b53c-8f3f42154362"
$ws.Range("S518").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S519:S520").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S521").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S522").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S523:S528").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S529").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S530").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S531:S532").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S533:S535").Value = "de503c24-f17d-47a9-9a47-6f0a194f8c9c"
$ws.Range("S536").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S537:S539").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S540").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S541:S551").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S552").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S553").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S554").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S555").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S556:S558").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S559:S561").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S562").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S563").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S564").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S565:S566").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S567:S568").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S569:S570").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S571:S572").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S573:S574").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S575:S576").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S577").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S578").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S579:S583").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S584:S586").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S587").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S588:S593").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S594:S595").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S596").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S597").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S598").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S599").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S600").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S601:S602").Value = "6afc31f0-3916-443a-92c4-b5eb425a9bc3"
$ws.Range("S603:S604").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S605:S606").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S607").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S608:S613").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S614").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S615:S616").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S617:S620").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S621:S622").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S623:S628").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S629:S631").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S632").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S633:S634").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S635:S636").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S637:S638").Value = "630f61e8-543f-46e2-af63-2b62e8bc4fd2"
$ws.Range("S639").Value = "de503c24-f17d-47a9-9a47-6f0a194f8c9c"
$ws.Range("S640").Value = "718c6b8f-7c00-4bcb-b53c-8f3f42154362"
$ws.Range("S641:S643").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S644:S645").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S646").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S647:S649").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S650").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S651:S652").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S653").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S654").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S655:S656").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S657:S658").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S659:S661").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S662").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S663").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S664").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S665:S667").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S668:S669").Value = "e0228b4f-7807-45db-a3f6-8c6e1f4adf41"
$ws.Range("S670:S671").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S672").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S673:S676").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S677").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S678").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S679").Value = "9cba5ad7-314e-4f2a-80a9-fc31cbf3f0c7"
$ws.Range("S680:S684").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S685:S687").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S688:S692").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S693").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S694").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S695:S696").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S697").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S698:S701").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S702").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S703:S704").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S705:S706").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S707").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S708:S712").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S713").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S714:S717").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S718:S719").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S720").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S721:S722").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S723:S734").Value = "7ba6e30b-04d4-4449-b7a9-2a6c7bb23764"
$ws.Range("S735").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S736:S737").Value = "4279fd55-c2c1-440d-abaa-430f3c27be44"
$ws.Range("S738").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S739").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S740").Value = "e42288a3-b5af-4464-bc45-85d438bcea11"
$ws.Range("S741:S743").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S744").Value = "e5ca6e2b-5f54-4acd-ad7b-03e631313986"
$ws.Range("S745").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S746:S747").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S748:S749").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S750:S751").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
$ws.Range("S752:S755").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S756:S758").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S759").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S760").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S761:S765").Value = "e5c40d19-b03a-4f5a-82c8-25540cd45e07"
$ws.Range("S766").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S767:S768").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S769").Value = "718c6b8f-7c00-4bcb-b53c-8f3f42154362"
$ws.Range("S770").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S771").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S772:S773").Value = "de503c24-f17d-47a9-9a47-6f0a194f8c9c"
$ws.Range("S774:S775").Value = "850a92da-c3d6-4fb9-a510-99626e9ad312"
$ws.Range("S776:S780").Value = "0eeb011c-24fb-4476-91f7-d8e28ae49c2f"
$ws.Range("S781").Value = "61f0969e-22a4-4374-8588-d6511915b05e"
$ws.Range("S782").Value = "5e1a20f6-82bf-4dee-aa79-41702d9feb41"
$ws.Range("S783").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S784").Value = "57ebee16-96d2-46a6-ab16-2476b305fd91"
$ws.Range("S785").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S786").Value = "7fc75193-58a0-4e7d-ab42-382ec10a8be4"
$ws.Range("S787").Value = "94742748-e7ab-454b-8ff6-9893440bd059"
$ws.Range("S788").Value = "3b8adc57-0f6f-482c-8306-9830e819d666"
$ws.Range("S789:S790").Value = "0435a227-38e8-494e-b1bf-271b00893eae"
$ws.Range("S791:S792").Value = "5184566d-523a-4432-848d-ac234ffb6ac6"
$ws.Range("S793").Value = "1ad987f9-6aa5-4e4d-8f1a-e8bea8fa4fcd"
$ws.Range("S794:S795").Value = "e7bb31c1-e095-453b-95ff-565ea62efb0a"
